# Weekly fruit/vegetable price update: insert two new price rows at the
# top of the date-ordered block (rows 1081-1082), pushing the existing
# rows 1081:1180 down to 1083:1182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 1081 (each Insert() pushes everything
# at/below that row down by one).
$ws.Rows.Item(1081).Insert()
$ws.Rows.Item(1081).Insert()

# Fill in the first new row (1081) with the new weekly record.
$ws.Cells.Item(1081, 1).Value = 10
$ws.Cells.Item(1081, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1081, 3).Value = "La Araucanía"
$ws.Cells.Item(1081, 4).Value = 45194
$ws.Cells.Item(1081, 5).Value = 9
$ws.Cells.Item(1081, 6).Value = 100112006
$ws.Cells.Item(1081, 7).Value = "Repollo"
$ws.Cells.Item(1081, 8).Value = "Crespo record"
$ws.Cells.Item(1081, 9).Value = "Primera"
$ws.Cells.Item(1081, 10).Value = 2000
$ws.Cells.Item(1081, 11).Value = 1200
$ws.Cells.Item(1081, 12).Value = 1200
$ws.Cells.Item(1081, 13).Value = 1200
$ws.Cells.Item(1081, 14).Value = "$/unidad"
$ws.Cells.Item(1081, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1081, 16).Value = 1200
$ws.Cells.Item(1081, 17).Value = 1
$ws.Cells.Item(1081, 18).Value = "Hortaliza"

# Fill in the second new row (1082) with the new weekly record.
$ws.Cells.Item(1082, 1).Value = 10
$ws.Cells.Item(1082, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1082, 3).Value = "La Araucanía"
$ws.Cells.Item(1082, 4).Value = 45194
$ws.Cells.Item(1082, 5).Value = 9
$ws.Cells.Item(1082, 6).Value = 100112006
$ws.Cells.Item(1082, 7).Value = "Repollo"
$ws.Cells.Item(1082, 8).Value = "Crespo record"
$ws.Cells.Item(1082, 9).Value = "Primera"
$ws.Cells.Item(1082, 10).Value = 3000
$ws.Cells.Item(1082, 11).Value = 1200
$ws.Cells.Item(1082, 12).Value = 1300
$ws.Cells.Item(1082, 13).Value = 1250
$ws.Cells.Item(1082, 14).Value = "$/unidad"
$ws.Cells.Item(1082, 15).Value = "Región del Maule"
$ws.Cells.Item(1082, 16).Value = 1250
$ws.Cells.Item(1082, 17).Value = 1
$ws.Cells.Item(1082, 18).Value = "Hortaliza"
